$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Col5a3", "Sdc3", "ECs", 3, 1, 3.389621, 10.168863, 0.03896472133961391, 0.03896472133961391, 3, 1, 34.36078833333333, 103.082365, 0.28490270239021, 0.28490270239021, 116.4700497112217, 1048.230447400995, 0.01110115440753749, 0.01110115440753748),
  @("ECs", "Col5a3", "Sdc3", "FAPs", 3, 1, 3.389621, 10.168863, 0.03896472133961391, 0.03896472133961391, 3, 1, 21.54461566666667, 64.633847, 0.17863732245739, 0.1786373224573899, 73.02808170066234, 657.252735305961, 0.006960553490406955, 0.006960553490406952),
  @("ECs", "Col5a3", "Sdc3", "M2", 3, 1, 3.389621, 10.168863, 0.03896472133961391, 0.03896472133961391, 3, 1, 60.03138866666666, 180.094166, 0.4977506538398792, 0.4977506538398792, 203.4836556836953, 1831.352901153258, 0.01939471552348151, 0.01939471552348151),
  @("ECs", "Col5a3", "Sdc3", "sCs", 3, 1, 3.389621, 10.168863, 0.03896472133961391, 0.03896472133961391, 3, 1, 4.668551, 14.005653, 0.03870932131252084, 0.03870932131252084, 15.824618509171, 142.421566582539, 0.001508297918187952, 0.001508297918187952),
  @("FAPs", "Col5a3", "Sdc3", "ECs", 3, 1, 71.87296166666665, 215.618885, 0.8262014907254878, 0.8262014907254879, 3, 1, 34.36078833333333, 103.082365, 0.28490270239021, 0.28490270239021, 2469.611622718113, 22226.50460446302, 0.2353870374265115, 0.2353870374265115),
  @("FAPs", "Col5a3", "Sdc3", "FAPs", 3, 1, 71.87296166666665, 215.618885, 0.8262014907254878, 0.8262014907254879, 3, 1, 21.54461566666667, 64.633847, 0.17863732245739, 0.1786373224573899, 1548.475335933399, 13936.27802340059, 0.1475904221135053, 0.1475904221135053),
  @("FAPs", "Col5a3", "Sdc3", "M2", 3, 1, 71.87296166666665, 215.618885, 0.8262014907254878, 0.8262014907254879, 3, 1, 60.03138866666666, 180.094166, 0.4977506538398792, 0.4977506538398792, 4314.6336964361, 38831.70326792491, 0.4112423322120944, 0.4112423322120945),
  @("FAPs", "Col5a3", "Sdc3", "sCs", 3, 1, 71.87296166666665, 215.618885, 0.8262014907254878, 0.8262014907254879, 3, 1, 4.668551, 14.005653, 0.03870932131252084, 0.03870932131252084, 335.5425870618783, 3019.883283556905, 0.03198169897337661, 0.03198169897337662),
  @("M2", "Col5a3", "Sdc3", "ECs", 2, 0.6666666666666666, 0.04149633333333334, 0.124489, 0.0004770129359444803, 0.0004770129359444803, 3, 1, 34.36078833333333, 103.082365, 0.28490270239021, 0.28490270239021, 1.425846726276111, 12.832620536485, 0.0001359022745256706, 0.0001359022745256706),
  @("M2", "Col5a3", "Sdc3", "FAPs", 2, 0.6666666666666666, 0.04149633333333334, 0.124489, 0.0004770129359444803, 0.0004770129359444803, 3, 1, 21.54461566666667, 64.633847, 0.17863732245739, 0.1786373224573899, 0.8940225532425558, 8.046202979183, 0.00008521231365466045, 0.00008521231365466043),
  @("M2", "Col5a3", "Sdc3", "M2", 2, 0.6666666666666666, 0.04149633333333334, 0.124489, 0.0004770129359444803, 0.0004770129359444803, 3, 1, 60.03138866666666, 180.094166, 0.4977506538398792, 0.4977506538398792, 2.491082514574889, 22.419742631174, 0.0002374335007564455, 0.0002374335007564455),
  @("M2", "Col5a3", "Sdc3", "sCs", 2, 0.6666666666666666, 0.04149633333333334, 0.124489, 0.0004770129359444803, 0.0004770129359444803, 3, 1, 4.668551, 14.005653, 0.03870932131252084, 0.03870932131252084, 0.1937277484796667, 1.743549736317, 0.00001846484700770381, 0.00001846484700770381),
  @("sCs", "Col5a3", "Sdc3", "ECs", 3, 1, 11.68797133333333, 35.063914, 0.1343567749989538, 0.1343567749989538, 3, 1, 34.36078833333333, 103.082365, 0.28490270239021, 0.28490270239021, 401.6079090307344, 3614.47118127661, 0.03827860828163535, 0.03827860828163535),
  @("sCs", "Col5a3", "Sdc3", "FAPs", 3, 1, 11.68797133333333, 35.063914, 0.1343567749989538, 0.1343567749989538, 3, 1, 21.54461566666667, 64.633847, 0.17863732245739, 0.1786373224573899, 251.8128502996842, 2266.315652697158, 0.02400113453982311, 0.02400113453982311),
  @("sCs", "Col5a3", "Sdc3", "M2", 3, 1, 11.68797133333333, 35.063914, 0.1343567749989538, 0.1343567749989538, 3, 1, 60.03138866666666, 180.094166, 0.4977506538398792, 0.4977506538398792, 701.6451498361914, 6314.806348525723, 0.06687617260354681, 0.06687617260354681),
  @("sCs", "Col5a3", "Sdc3", "sCs", 3, 1, 11.68797133333333, 35.063914, 0.1343567749989538, 0.1343567749989538, 3, 1, 4.668551, 14.005653, 0.03870932131252084, 0.03870932131252084, 54.56589025620466, 491.093012305842, 0.00520085957394857, 0.005200859573948571)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
